$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow writing to locked cells, re-protect after.
$ws.Unprotect()

# Update the confidential disclaimer date (2021-04-21 -> 2021-04-22) in A42.
$ws.Range("A42").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-22 for illustrative purposes only and are subject to change."

# Updated Weight (D) and Percent Change (E) values for holdings rows 2-39 (row 39 = Total).
$ws.Range("D2").Value = 0.06246348000091403
$ws.Range("E2").Value = -0.01168539325842699
$ws.Range("D3").Value = 0.05670841779858262
$ws.Range("E3").Value = -0.01308619234016406
$ws.Range("D4").Value = 0.2901252098709818
$ws.Range("E4").Value = -0.01706484641638228
$ws.Range("D5").Value = 0.03658278356113109
$ws.Range("E5").Value = -0.01575838335286528
$ws.Range("D6").Value = 0.03217405118730122
$ws.Range("E6").Value = 0.0005275884809847753
$ws.Range("D7").Value = 0.02948498232350436
$ws.Range("E7").Value = -0.02105752623887336
$ws.Range("D8").Value = 0.02900316316237894
$ws.Range("E8").Value = -0.008463893390959876
$ws.Range("D9").Value = 0.02458278791360769
$ws.Range("E9").Value = -0.0108356940509915
$ws.Range("D10").Value = 0.02479116273148375
$ws.Range("E10").Value = -0.01133715188623341
$ws.Range("D11").Value = 0.02296246968228009
$ws.Range("E11").Value = -0.01641954423325709
$ws.Range("D12").Value = 0.02190862628221241
$ws.Range("E12").Value = -0.009297520661157077
$ws.Range("D13").Value = 0.02191700479969359
$ws.Range("E13").Value = -0.00191141936540884
$ws.Range("D14").Value = 0.02233636592141386
$ws.Range("E14").Value = -0.01622214103032527
$ws.Range("D15").Value = 0.02128034628303938
$ws.Range("E15").Value = -0.01377511888326433
$ws.Range("D16").Value = 0.02169763997836819
$ws.Range("E16").Value = -0.0056919335021689
$ws.Range("D17").Value = 0.02084292238337267
$ws.Range("E17").Value = 0.0006786739754633953
$ws.Range("D18").Value = 0.01663516561717575
$ws.Range("E18").Value = -0.0177394034536893
$ws.Range("D19").Value = 0.0166519226521381
$ws.Range("E19").Value = -0.004665629860031162
$ws.Range("D20").Value = 0.01581581189466572
$ws.Range("E20").Value = -0.01479188166494672
$ws.Range("D21").Value = 0.01584301487350071
$ws.Range("E21").Value = -0.01303571428571426
$ws.Range("D22").Value = 0.0161938244885568
$ws.Range("E22").Value = -0.0328307262269526
$ws.Range("D23").Value = 0.01544976861146203
$ws.Range("E23").Value = -0.003112982970151834
$ws.Range("D24").Value = 0.01408820511481289
$ws.Range("E24").Value = 0.04151444702756568
$ws.Range("D25").Value = 0.01412650690901257
$ws.Range("E25").Value = -0.008280377431157282
$ws.Range("D26").Value = 0.01467437490274935
$ws.Range("E26").Value = 0.007340946166394913
$ws.Range("D27").Value = 0.01272120102239676
$ws.Range("E27").Value = -0.005542725173210195
$ws.Range("D28").Value = 0.01341389767545106
$ws.Range("E28").Value = -0.01732697362016933
$ws.Range("D29").Value = 0.01462453904552364
$ws.Range("E29").Value = -0.02251454591449531
$ws.Range("D30").Value = 0.01334643428794027
$ws.Range("E30").Value = -0.02334985650926158
$ws.Range("D31").Value = 0.0127945402533359
$ws.Range("E31").Value = -0.006191318546741043
$ws.Range("D32").Value = 0.01326199624163645
$ws.Range("E32").Value = -0.006768953068592043
$ws.Range("D33").Value = 0.01293175207857962
$ws.Range("E33").Value = -0.01602928183768793
$ws.Range("D34").Value = 0.006685621702318892
$ws.Range("E34").Value = -0.03321831971615508
$ws.Range("D35").Value = 0.00553743837165145
$ws.Range("E35").Value = -0.0002358027117311812
$ws.Range("D36").Value = 0.005764637650881323
$ws.Range("E36").Value = -0.009890898108648849
$ws.Range("D37").Value = 0.005574869670528403
$ws.Range("E37").Value = -0.009856735761408442
$ws.Range("D38").Value = 0.005003063055416822
$ws.Range("E38").Value = -0.01017855977728965
$ws.Range("D39").Value = 1
$ws.Range("E39").Value = -0.01230749812027399

# Restore sheet protection (best effort; original password cannot be recovered from its stored hash).
$ws.Protect()

